$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Layout")
Write-Output $ws.Range("B5").Value
